# The workbook tracks daily working hours. A new day-entry (2014-07-29,
# 18:40 -> 20:30) needs to be inserted right before the blank separator row
# that precedes the summary block, pushing the separator and the three
# summary rows (sum[min]/sum[h]/sum[working weeks]) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 166 - this shifts the old row 166 (blank separator)
# and rows 167-169 (the summary rows) down to 167-170, carrying their
# formatting/formulas with them (Excel auto-adjusts the SUM ranges that
# don't cross the insertion point).
$ws.Rows(166).Insert()

# Fill in the new data row with the new reading.
$ws.Range("A166").Value = 2014
$ws.Range("B166").Value = 7
$ws.Range("C166").Value = 29
$ws.Range("D166").Value = 0.77777777777777779
$ws.Range("E166").Value = 0.85416666666666663
$ws.Range("F166").Formula = "=(E166-D166)*24*60"
$ws.Range("G166").Formula = "=F166/60"

# The "sum [min]" total (now on row 168) must include the new row 166, so
# extend its range from F2:F165 to F2:F166.
$ws.Range("F168").Formula = "=SUM(F2:F166)"

# Restore the view: select the new blank separator's E cell and scroll down
# a couple of rows, matching where the user was working.
[void]$ws.Range("E167").Select()
$excel.ActiveWindow.ScrollRow = 147
